$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# idx16: '2016-08-13 22:20:38' -> '2016-08-13 22:21:46'
$overview.Range("G4").Value = "2016-08-13 22:21:46"
$overview.Range("G5").Value = "2016-08-13 22:21:46"
$dede.Range("H4").Value = "2016-08-13 22:21:46"
$dede.Range("H5").Value = "2016-08-13 22:21:46"

# idx43: 'ht' -> 'mt'
$zhcn.Range("E4").Value = "mt"
$zhcn.Range("E5").Value = "mt"
$dede.Range("E4").Value = "mt"
$dede.Range("E5").Value = "mt"

# idx45: '2016-08-13 22:20:30' -> '2016-08-13 22:21:37'
$zhcn.Range("H4").Value = "2016-08-13 22:21:37"
$zhcn.Range("H5").Value = "2016-08-13 22:21:37"

# idx46: '2016-08-13 22:20:58' -> '2016-08-13 22:22:10'
$zhcn.Range("K4").Value = "2016-08-13 22:22:10"
$zhcn.Range("K5").Value = "2016-08-13 22:22:10"

# idx52: '2016-08-13 22:21:12' -> '2016-08-13 22:22:20'
$dede.Range("K4").Value = "2016-08-13 22:22:20"
$dede.Range("K5").Value = "2016-08-13 22:22:20"
